# Rename the header labels in row 1 (B1:D1): " LN1"/" LN2"/" LN3" -> "Node_1"/"Node_2"/"Node_3"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Node_1"
$ws.Range("C1").Value = "Node_2"
$ws.Range("D1").Value = "Node_3"

# Remove the trailing, empty formatted row (row 4 only contained a styled,
# valueless cell B4) so the sheet's used range shrinks back to A1:D3.
$ws.Rows.Item(4).Delete()

# Match the saved selection/active cell from the edit (D1).
[void]$ws.Range("D1").Select()
